$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a brand-new row at position 28 (shifts old rows 28-34 down to 29-35,
# matching the defined-name / dimension growing from R34 to R35).
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the "nicardipine" entry.
$ws.Range("A28").Value = "nicardipine"
$ws.Range("B28").Value = "mg"
$ws.Range("C28").Value = "microg/kg/min"
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 50
$ws.Range("H28").Value = 10
$ws.Range("I28").Value = 50
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 50
$ws.Range("L28").Value = 0.5
$ws.Range("M28").Value = 3
$ws.Range("N28").Value = 5
$ws.Range("P28").Value = 0.2
$ws.Range("Q28").Value = 0

# Move/restore the visible selection to match the saved view state.
$ws.Range("D25").Select()

# The admin table's defined name covered A4:R34; grow it by one row to A4:R35
# now that a new medication row has been inserted.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*Tbl_Admin_PedMedCont*") {
        $n.RefersTo = "=Tbl_Admin_PedMedCont!`$A`$4:`$R`$35"
    }
}
